$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "总计" (Total) sheet: insert a new row 2 for 2022-Q3 and shift the rest
#    of the quarters down by one row.
# ---------------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item(1)

$wsTotal.Rows.Item(2).Insert()

# Give the new row 2 the same look as the other data rows: column A uses the
# bold/centered/bordered style, columns B:D are plain.
$wsTotal.Range("B2:D2").ClearFormats()
$wsTotal.Range("A3").Copy()
$wsTotal.Range("A2").PasteSpecial(-4122)   # xlPasteFormats

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q3"
$wsTotal.Range("C2").Value = 6
$wsTotal.Range("D2").Value = 1.17

# The row-index column (A) holds each quarter's 0-based position; after the
# insert the rows below still carry their old (pre-shift) index, so bump
# them by one to keep the 0..7 sequence intact.
$wsTotal.Range("A3").Value = 1
$wsTotal.Range("A4").Value = 2
$wsTotal.Range("A5").Value = 3
$wsTotal.Range("A6").Value = 4
$wsTotal.Range("A7").Value = 5
$wsTotal.Range("A8").Value = 6
$wsTotal.Range("A9").Value = 7

# ---------------------------------------------------------------------------
# 2. Create the new "2022-Q3" sheet right after "总计" by duplicating the
#    "2022-Q2" sheet (this keeps column/row styling identical to the other
#    quarter sheets) and then overwrite its contents with the Q3 data.
# ---------------------------------------------------------------------------
$wsQ2 = $wb.Worksheets.Item("2022-Q2")
$wsQ2.Copy($null, $wsTotal)

$wsQ3 = $wb.Worksheets.Item(2)
$wsQ3.Name = "2022-Q3"

# Drop the extra fund rows coming from the copied 2022-Q2 sheet (17 rows ->
# 7 rows: 1 header + 6 funds) before writing the new values.
$wsQ3.Rows("8:17").Delete()

# Columns that must keep their values as text (to match the original text
# formatted numbers such as "007777" or "0.0270").
$wsQ3.Range("B2:B7").NumberFormat = "@"
$wsQ3.Range("D2:D7").NumberFormat = "@"
$wsQ3.Range("E2:E7").NumberFormat = "@"
$wsQ3.Range("F2:F7").NumberFormat = "@"
$wsQ3.Range("G2:G6").NumberFormat = "@"

$wsQ3.Range("B2").Value = "007777"
$wsQ3.Range("C2").Value = "中邮研究精选混合"
$wsQ3.Range("D2").Value = "14.89"
$wsQ3.Range("E2").Value = "70.15"
$wsQ3.Range("F2").Value = "4.25"
$wsQ3.Range("G2").Value = "0.6328"
$wsQ3.Range("H2").Value = 5

$wsQ3.Range("B3").Value = "012719"
$wsQ3.Range("C3").Value = "华夏新兴经济一年持有混合A"
$wsQ3.Range("D3").Value = "17.53"
$wsQ3.Range("E3").Value = "87.21"
$wsQ3.Range("F3").Value = "2.76"
$wsQ3.Range("G3").Value = "0.4838"
$wsQ3.Range("H3").Value = 6

$wsQ3.Range("B4").Value = "012720"
$wsQ3.Range("C4").Value = "华夏新兴经济一年持有混合C"
$wsQ3.Range("D4").Value = "0.98"
$wsQ3.Range("E4").Value = "87.21"
$wsQ3.Range("F4").Value = "2.76"
$wsQ3.Range("G4").Value = "0.0270"
$wsQ3.Range("H4").Value = 6

$wsQ3.Range("B5").Value = "001266"
$wsQ3.Range("C5").Value = "国投瑞银招财灵活配置混合"
$wsQ3.Range("D5").Value = "0.69"
$wsQ3.Range("E5").Value = "68.03"
$wsQ3.Range("F5").Value = "2.62"
$wsQ3.Range("G5").Value = "0.0181"
$wsQ3.Range("H5").Value = 4

$wsQ3.Range("B6").Value = "002409"
$wsQ3.Range("C6").Value = "华夏新活力灵活配置混合A"
$wsQ3.Range("D6").Value = "0.15"
$wsQ3.Range("E6").Value = "69.89"
$wsQ3.Range("F6").Value = "3.11"
$wsQ3.Range("G6").Value = "0.0047"
$wsQ3.Range("H6").Value = 9

$wsQ3.Range("B7").Value = "002410"
$wsQ3.Range("C7").Value = "华夏新活力灵活配置混合C"
$wsQ3.Range("D7").Value = "0.00"
$wsQ3.Range("E7").Value = "69.89"
$wsQ3.Range("F7").Value = "3.11"
$wsQ3.Range("G7").Value = 0
$wsQ3.Range("H7").Value = 9
